$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 100
$ws.Range("I42").Value = 100
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 300
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -70
$ws.Range("H107").Value = 1174.9
$ws.Range("I107").Value = 1174.9
$ws.Range("K107").Value = 1174.9
$ws.Range("M107").Value = 745.0999999999999
$ws.Range("H113").Value = 2894.2334
$ws.Range("I113").Value = 2626.353
$ws.Range("J113").Value = 3244.5386
$ws.Range("K113").Value = 2626.353
$ws.Range("L113").Value = 3244.5386
$ws.Range("M113").Value = 627.6469999999999
$ws.Range("N113").Value = -9752.5386
$ws.Range("H135").Value = 28177.37
$ws.Range("I135").Value = 37333.965
$ws.Range("J135").Value = 2538.9
$ws.Range("K135").Value = 336005.6849999999
$ws.Range("L135").Value = 22850.1
$ws.Range("M135").Value = -333470.6849999999
$ws.Range("N135").Value = -27920.1
$ws.Range("H137").Value = 1854362.8
$ws.Range("I137").Value = 3126687
$ws.Range("K137").Value = 9380061
$ws.Range("M137").Value = -9377511
$ws.Range("H138").Value = 3177067.5
$ws.Range("I138").Value = 4222.75
$ws.Range("J138").Value = 3392175.5
$ws.Range("K138").Value = 12668.25
$ws.Range("L138").Value = 10176526.5
$ws.Range("M138").Value = -7528.25
$ws.Range("N138").Value = -10186806.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 10563
$ws.Range("I21").Value = 8253.75
$ws.Range("J21").Value = 19800
$ws.Range("K21").Value = 8253.75
$ws.Range("L21").Value = 19800
$ws.Range("M21").Value = -7879.75
$ws.Range("N21").Value = -20548
$ws.Range("H133").Value = 49796.668
$ws.Range("J133").Value = 49796.668
$ws.Range("L133").Value = 49796.668
$ws.Range("N133").Value = -54856.668
$ws.Range("H139").Value = 49670.715
$ws.Range("J139").Value = 49670.715
$ws.Range("L139").Value = 49670.715
$ws.Range("N139").Value = -59950.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 318.25925
$ws.Range("J80").Value = 357.35
$ws.Range("L80").Value = 357.35
$ws.Range("N80").Value = -2353.35
$ws.Range("H83").Value = 318.25925
$ws.Range("J83").Value = 357.35
$ws.Range("L83").Value = 1786.75
$ws.Range("N83").Value = -11770.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1349.9231
$ws.Range("I16").Value = 1219.9
$ws.Range("J16").Value = 1783.3334
$ws.Range("K16").Value = 1219.9
$ws.Range("L16").Value = 1783.3334
$ws.Range("M16").Value = -932.9000000000001
$ws.Range("N16").Value = -2357.3334
$ws.Range("H52").Value = 48900
$ws.Range("J52").Value = 48900
$ws.Range("L52").Value = 48900
$ws.Range("N52").Value = -49488
$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 298
$ws.Range("H113").Value = 1349.9231
$ws.Range("I113").Value = 1219.9
$ws.Range("J113").Value = 1783.3334
$ws.Range("K113").Value = 1219.9
$ws.Range("L113").Value = 1783.3334
$ws.Range("M113").Value = 950.0999999999999
$ws.Range("N113").Value = -6123.3334
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -1130
$ws.Range("H141").Value = 45862.715
$ws.Range("J141").Value = 45862.715
$ws.Range("L141").Value = 45862.715
$ws.Range("N141").Value = -56222.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 568.625
$ws.Range("J127").Value = 568.625
$ws.Range("L127").Value = 1705.875
$ws.Range("N127").Value = -11625.875
$ws.Range("H137").Value = 1684.6842
$ws.Range("J137").Value = 2410
$ws.Range("L137").Value = 7230
$ws.Range("N137").Value = -17430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 31666.666
$ws.Range("I46").Value = 15000
$ws.Range("J46").Value = 40000
$ws.Range("K46").Value = 15000
$ws.Range("L46").Value = 40000
$ws.Range("M46").Value = -14844
$ws.Range("N46").Value = -40312
$ws.Range("H54").Value = 4886.5
$ws.Range("I54").Value = 300
$ws.Range("K54").Value = 300
$ws.Range("M54").Value = 90
$ws.Range("H102").Value = 2320.6428
$ws.Range("I102").Value = 2297.9167
$ws.Range("J102").Value = 2457
$ws.Range("K102").Value = 2297.9167
$ws.Range("L102").Value = 2457
$ws.Range("M102").Value = -675.9167000000002
$ws.Range("N102").Value = -5701
$ws.Range("H122").Value = 2631.077
$ws.Range("I122").Value = 2382.182
$ws.Range("K122").Value = 7146.545999999999
$ws.Range("M122").Value = -4696.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 642.13336
$ws.Range("I22").Value = 520.73334
$ws.Range("J22").Value = 763.5333000000001
$ws.Range("K22").Value = 520.73334
$ws.Range("L22").Value = 763.5333000000001
$ws.Range("M22").Value = -225.73334
$ws.Range("N22").Value = -1353.5333
$ws.Range("H27").Value = 642.13336
$ws.Range("I27").Value = 520.73334
$ws.Range("J27").Value = 763.5333000000001
$ws.Range("K27").Value = 520.73334
$ws.Range("L27").Value = 763.5333000000001
$ws.Range("M27").Value = -413.73334
$ws.Range("N27").Value = -977.5333000000001
$ws.Range("H46").Value = 3788339
$ws.Range("I46").Value = 6061082
$ws.Range("J46").Value = 433.33334
$ws.Range("K46").Value = 6061082
$ws.Range("L46").Value = 433.33334
$ws.Range("M46").Value = -6060894
$ws.Range("N46").Value = -809.33334
$ws.Range("H132").Value = 85616.234
$ws.Range("I132").Value = 7750.875
$ws.Range("J132").Value = 210200.8
$ws.Range("K132").Value = 23252.625
$ws.Range("L132").Value = 630602.3999999999
$ws.Range("M132").Value = -20722.625
$ws.Range("N132").Value = -635662.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13993.2
$ws.Range("J54").Value = 13993.2
$ws.Range("L54").Value = 13993.2
$ws.Range("N54").Value = -15033.2
$ws.Range("H81").Value = 2034.3704
$ws.Range("I81").Value = 1100
$ws.Range("J81").Value = 2246.7273
$ws.Range("K81").Value = 2200
$ws.Range("L81").Value = 4493.4546
$ws.Range("M81").Value = -1139
$ws.Range("N81").Value = -6615.4546
$ws.Range("H84").Value = 2034.3704
$ws.Range("I84").Value = 1100
$ws.Range("J84").Value = 2246.7273
$ws.Range("K84").Value = 11000
$ws.Range("L84").Value = 22467.273
$ws.Range("M84").Value = -5696
$ws.Range("N84").Value = -33075.273
$ws.Range("H126").Value = 1001.35
$ws.Range("I126").Value = 565
$ws.Range("J126").Value = 1655.875
$ws.Range("K126").Value = 1695
$ws.Range("L126").Value = 4967.625
$ws.Range("M126").Value = 775
$ws.Range("N126").Value = -9907.625
